# Version Feb 5, 2022
# Applies the diff: updates Trial!row3 values/formats, appends two new
# data rows (4 and 5) to the Trial sheet with the same layout/styling as
# row 3, fixes the two boolean cells on RFC_params to be TRUE()/FALSE()
# formulas, and restores the recorded cell selections on both sheets.

$wb = $excel.ActiveWorkbook

$trial = $wb.Worksheets.Item("Trial")
$rfc   = $wb.Worksheets.Item("RFC_params")

# ---------------------------------------------------------------------
# Trial sheet: row 3 updates (values + number formats)
# ---------------------------------------------------------------------
$trial.Range("D3").NumberFormat = "0.00"
$trial.Range("D3").Value = 0.4

$trial.Range("E3").Value = 0.77

$trial.Range("G3").NumberFormat = "0.00"
$trial.Range("G3").Value = 0.4

$trial.Range("H3").Value = 0.77

$trial.Range("I3").NumberFormat = "0.00"
$trial.Range("I3").Value = 0.6

$trial.Range("J3").Value = 0.001493

$trial.Range("K3").Value = 0.9

$trial.Range("L3").Value = 0.001789

# ---------------------------------------------------------------------
# Trial sheet: brand new rows 4 and 5, cloned (format-only) from the now
# up-to-date row 3 so every cell keeps the same style index, then filled
# in with their own values.
# ---------------------------------------------------------------------
$trial.Range("A3:L3").Copy()
$trial.Range("A4").PasteSpecial(-4122)
$trial.Range("A5").PasteSpecial(-4122)

$trial.Rows.Item(4).RowHeight = 35.05
$trial.Rows.Item(5).RowHeight = 35.05

# Row 4
$trial.Range("A4").Value = 1
$trial.Range("B4").Value = 2
$trial.Range("C4").Value = "RFC"
$trial.Range("D4").Value = 0.45
$trial.Range("E4").Value = 0.8
$trial.Range("F4").Value = "over-fitting (high variance)"
$trial.Range("G4").Value = 0.45
$trial.Range("H4").Value = 0.8
$trial.Range("I4").Value = 0.55
$trial.Range("J4").Value = 0.001342
$trial.Range("K4").Value = 0.9
$trial.Range("L4").Value = 0.000575

# Row 5
$trial.Range("A5").Value = 1
$trial.Range("B5").Value = 3
$trial.Range("C5").Value = "RFC"
$trial.Range("D5").Value = 0.45
$trial.Range("E5").Value = 0.8
$trial.Range("F5").Value = "over-fitting (high variance)"
$trial.Range("G5").Value = 0.45
$trial.Range("H5").Value = 0.8
$trial.Range("I5").Value = 0.55
$trial.Range("J5").Value = 0.001138
$trial.Range("K5").Value = 0.9
$trial.Range("L5").Value = 0.000831

# ---------------------------------------------------------------------
# RFC_params sheet: B3 / O3 become TRUE()/FALSE() formulas
# ---------------------------------------------------------------------
$rfc.Range("B3").Formula = "=TRUE()"
$rfc.Range("O3").Formula = "=FALSE()"

# ---------------------------------------------------------------------
# RFC_params sheet: page setup cosmetic tweak (firstPageNumber 1 -> 0)
# ---------------------------------------------------------------------
$rfc.PageSetup.FirstPageNumber = 0

# ---------------------------------------------------------------------
# Selections: Trial -> J4, RFC_params -> B3 (Trial stays the active tab)
# ---------------------------------------------------------------------
$rfc.Activate()
$rfc.Range("B3").Select()

$trial.Activate()
$trial.Range("J4").Select()
